$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Work on the existing "ODI Batting" sheet: rename MATCH_CARD_LINK header
#    to MATCH_CODE, replace the howstat URL with the bare numeric match code,
#    and drop the (already blank) INNING_NUMBER cells that have no value.
# ---------------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")

$battingLastRow = $batting.Cells.Item(1, 1).CurrentRegion.Rows.Count

$batting.Range("D1").Value = "MATCH_CODE"

$batting.Range("D2:D" + $battingLastRow).NumberFormat = "@"
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $batting.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val -match "MatchCode=(\d+)") {
        $cell.Value = $matches[1]
    }

    $bcell = $batting.Cells.Item($r, 2)
    $bval = $bcell.Value2
    if ($bval -eq $null -or $bval -eq "") {
        $bcell.ClearContents()
    }
}

# ---------------------------------------------------------------------------
# 2. Work on the existing "ODI Bowling" sheet: rename MATCH_CARD_LINK header
#    to MATCH_CODE and replace the howstat URL with the bare numeric code.
# ---------------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")

$bowlingLastRow = $bowling.Cells.Item(1, 1).CurrentRegion.Rows.Count

$bowling.Range("B1").Value = "MATCH_CODE"

$bowling.Range("B2:B" + $bowlingLastRow).NumberFormat = "@"
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowling.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -match "MatchCode=(\d+)") {
        $cell.Value = $matches[1]
    }
}

# ---------------------------------------------------------------------------
# 3. Insert a new "Player Info" sheet before "ODI Batting" with player bio
#    details.
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($batting)
$playerInfo.Name = "Player Info"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "3149"
$playerInfo.Range("B2").Value = "Ravindranath Rampaul"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast Medium"

# ---------------------------------------------------------------------------
# 4. Append a new "ODI Batting Extra" sheet after "ODI Bowling" with extra
#    per-match batting detail.
# ---------------------------------------------------------------------------
$extra = $wb.Worksheets.Add()
$extra.Name = "ODI Batting Extra"
$extra.Move($null, $bowling)

$extraHeader = $extra.Range("A1:F1")
$extraHeader.Font.Bold = $true
$extraHeader.HorizontalAlignment = -4108
$extraHeader.VerticalAlignment = -4160
$extraHeader.Borders.LineStyle = 1

$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

$extra.Range("A2:A21").NumberFormat = "@"

$extraRows = @(
    @("3452", "", "", "", "", "NO"),
    @("3506", 10, "0", "0", "0.86%", "NO"),
    @("3510", "", "", "", "", "NO"),
    @("3513", "", "", "", "", "NO"),
    @("3520", "", "", "", "", "NO"),
    @("3580", 11, "", "", "", "NO"),
    @("3581", "", "", "", "", "NO"),
    @("3583", "", "", "", "", "NO"),
    @("3593", 11, "", "", "", "NO"),
    @("3596", 11, "0", "0", "0.63%", "NO"),
    @("3622", 10, "", "", "", "NO"),
    @("3625", 10, "", "", "", "NO"),
    @("3629", "", "", "", "", "NO"),
    @("3655", 10, "", "", "", "NO"),
    @("3657", 9, "1", "1", "7.44%", "NO"),
    @("3661", 11, "0", "0", "0.47%", "NO"),
    @("3678", 11, "0", "0", "0.97%", "NO"),
    @("3680", "", "", "", "", ""),
    @("3853", "", "", "", "", ""),
    @("3855", "", "", "", "", "")
)

$r = 2
foreach ($row in $extraRows) {
    $extra.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne "") {
        $extra.Cells.Item($r, 2).Value = $row[1]
    }
    if ($row[2] -ne "") {
        $extra.Cells.Item($r, 3).Value = $row[2]
    }
    if ($row[3] -ne "") {
        $extra.Cells.Item($r, 4).Value = $row[3]
    }
    if ($row[4] -ne "") {
        $extra.Cells.Item($r, 5).Value = $row[4]
    }
    if ($row[5] -ne "") {
        $extra.Cells.Item($r, 6).Value = $row[5]
    }
    $r = $r + 1
}
